$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# K2:K31 -> new "Ship To Customer Name" value (plain text, no format change needed)
$ws.Range("K2:K31").Value = "mnageq256303"

# AX2:AX31 -> new "Previous Doc" value. The source value is a purely numeric
# string, so a leading apostrophe (text/quote prefix) is required - exactly
# like typing '9824427933 into a cell in Excel - otherwise it is parsed back
# into a number, same as the original 2140999000.
$ws.Range("AX2:AX31").Value = "'9824427933"
